# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) in the EC (Estado de Cuenta) table
# is updated so the arrears periods are listed in ascending chronological
# order (2412, 2501, 2502, 2503, 2504, 2505, 2506) instead of the
# previous scrambled order (2505, 2504, 2503, 2502, 2501, 2412, 2506).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2412"
$ws.Range("E17").Value = "2501"
$ws.Range("E18").Value = "2502"
$ws.Range("E19").Value = "2503"
$ws.Range("E20").Value = "2504"
$ws.Range("E21").Value = "2505"
$ws.Range("E22").Value = "2506"
